$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.902.66"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = "'  -0.66%  "
$ws.Range('E2').ClearFormats()

$ws.Range('D3').Value = "'2.347.18"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = "'  -1.29%  "
$ws.Range('E3').ClearFormats()

$ws.Range('E4').Value = "'  -0.04%  "
$ws.Range('E4').ClearFormats()

$ws.Range('D5').Value = "'240.15"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = "'  -1.39%  "
$ws.Range('E5').ClearFormats()

$ws.Range('D6').Value = "'0.670"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = "'  -3.77%  "
$ws.Range('E6').ClearFormats()

$ws.Range('D7').Value = "'72.27"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = "'  -5.46%  "
$ws.Range('E7').ClearFormats()

$ws.Range('E8').Value = "'  -0.05%  "
$ws.Range('E8').ClearFormats()

$ws.Range('D9').Value = "'0.594"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = "'  -1.41%  "
$ws.Range('E9').ClearFormats()

$ws.Range('D10').Value = "'0.100"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = "'  -3.43%  "
$ws.Range('E10').ClearFormats()

$ws.Range('D11').Value = "'58.13"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = "'  +0.21%  "
$ws.Range('E11').ClearFormats()

$ws.Range('D12').Value = "'32.88"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = "'  +1.67%  "
$ws.Range('E12').ClearFormats()

$ws.Range('E13').Value = "'  -0.12%  "
$ws.Range('E13').ClearFormats()

$ws.Range('D14').Value = "'7.24"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = "'  -3.72%  "
$ws.Range('E14').ClearFormats()

$ws.Range('D15').Value = "'2.698.03"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = "'  -1.45%  "
$ws.Range('E15').ClearFormats()

$ws.Range('D16').Value = "'16.30"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = "'  -5.04%  "
$ws.Range('E16').ClearFormats()

$ws.Range('D17').Value = "'0.903"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = "'  -2.79%  "
$ws.Range('E17').ClearFormats()

$ws.Range('D18').Value = "'2.348.81"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = "'  -1.54%  "
$ws.Range('E18').ClearFormats()

$ws.Range('D19').Value = "'43.792.55"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = "'  -1.78%  "
$ws.Range('E19').ClearFormats()

$ws.Range('E20').Value = "'  -1.06%  "
$ws.Range('E20').ClearFormats()

$ws.Range('D21').Value = "'6.72"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = "'  -0.49%  "
$ws.Range('E21').ClearFormats()

$ws.Range('D22').Value = "'78.18"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "'  -0.97%  "
$ws.Range('E22').ClearFormats()

$ws.Range('D23').Value = "'254.46"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = "'  -1.31%  "
$ws.Range('E23').ClearFormats()

$ws.Range('D24').Value = "'1.93"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = "'  +6.75%  "
$ws.Range('E24').ClearFormats()

$ws.Range('E25').Value = "'  +0.00%  "
$ws.Range('E25').ClearFormats()

$ws.Range('E26').Value = "'  +0.49%  "
$ws.Range('E26').ClearFormats()

$ws.Range('D27').Value = "'2.49"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = "'  -3.55%  "
$ws.Range('E27').ClearFormats()

$ws.Range('D28').Value = "'10.46"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = "'  -4.25%  "
$ws.Range('E28').ClearFormats()

$ws.Range('D30').Value = "'176.82"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = "'  +0.75%  "
$ws.Range('E30').ClearFormats()

$ws.Range('D31').Value = "'22.40"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = "'  -3.60%  "
$ws.Range('E31').ClearFormats()

$ws.Range('D32').Value = "'0.127"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = "'  -3.10%  "
$ws.Range('E32').ClearFormats()

$ws.Range('E33').Value = "'  +0.28%  "
$ws.Range('E33').ClearFormats()

$ws.Range('D34').Value = "'0.0748"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = "'  -1.47%  "
$ws.Range('E34').ClearFormats()

$ws.Range('D35').Value = "'5.13"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = "'  -4.71%  "
$ws.Range('E35').ClearFormats()

$ws.Range('D36').Value = "'5.39"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = "'  +0.68%  "
$ws.Range('E36').ClearFormats()

$ws.Range('D37').Value = "'3.75"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = "'  -3.89%  "
$ws.Range('E37').ClearFormats()

$ws.Range('D38').Value = "'6.43"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = "'  -2.76%  "
$ws.Range('E38').ClearFormats()

$ws.Range('D39').Value = "'2.38"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = "'  -5.03%  "
$ws.Range('E39').ClearFormats()

$ws.Range('D40').Value = "'0.0275"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "'  -1.17%  "
$ws.Range('E40').ClearFormats()

$ws.Range('D41').Value = "'66.78"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = "'  +22.19%  "
$ws.Range('E41').ClearFormats()

$ws.Range('D42').Value = "'5.17"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = "'  +15.34%  "
$ws.Range('E42').ClearFormats()

$ws.Range('E43').Value = "'  +7.99%  "
$ws.Range('E43').ClearFormats()

$ws.Range('D44').Value = "'9.19"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "'  +0.33%  "
$ws.Range('E44').ClearFormats()

$ws.Range('D45').Value = "'18.81"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = "'  -1.39%  "
$ws.Range('E45').ClearFormats()

$ws.Range('D46').Value = "'0.199"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = "'  +1.86%  "
$ws.Range('E46').ClearFormats()

$ws.Range('D47').Value = "'2.48"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = "'  -2.33%  "
$ws.Range('E47').ClearFormats()

$ws.Range('B48').Value = "'BinanceUSD"
$ws.Range('B48').ClearFormats()
$ws.Range('C48').Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range('C48').ClearFormats()
$ws.Range('D48').Value = "'1.00"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = "'  -0.06%  "
$ws.Range('E48').ClearFormats()

$ws.Range('B49').Value = "'TrustWalletToken"
$ws.Range('B49').ClearFormats()
$ws.Range('C49').Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range('C49').ClearFormats()
$ws.Range('D49').Value = "'1.24"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = "'  -2.86%  "
$ws.Range('E49').ClearFormats()

$ws.Range('D50').Value = "'99.13"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = "'  -3.92%  "
$ws.Range('E50').ClearFormats()

$ws.Range('E51').Value = "'  -5.78%  "
$ws.Range('E51').ClearFormats()
